$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.056501879908353
$ws.Range("C2").Value = 0.1897083962827537
$ws.Range("E2").Value = 0.08946581615927029
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.8393329927350806
$ws.Range("H2").Value = 0.8837776229843968
$ws.Range("I2").Value = 0.8150015711907415
$ws.Range("L2").Value = 0.2075937874129465
$ws.Range("M2").Value = 0.2389369135504111
$ws.Range("N2").Value = 1.408836080745495
# Row 3
$ws.Range("B3").Value = 0.9703779529608028
$ws.Range("C3").Value = 0.171470773973823
$ws.Range("E3").Value = 0.08987119776263774
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.8341860943733934
$ws.Range("H3").Value = 0.8869222500388787
$ws.Range("I3").Value = 0.8209674896580594
$ws.Range("L3").Value = 0.2048676727839691
$ws.Range("M3").Value = 0.2253560573279287
$ws.Range("N3").Value = 1.425743839560418
# Row 4
$ws.Range("B4").Value = 0.9178159457725599
$ws.Range("C4").Value = 0.1601919505871479
$ws.Range("E4").Value = 0.09013542060389934
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.8317159349234657
$ws.Range("H4").Value = 0.8893691384776332
$ws.Range("I4").Value = 0.8251536703059053
$ws.Range("L4").Value = 0.2032991707167895
$ws.Range("M4").Value = 0.2171178492876109
$ws.Range("N4").Value = 1.436700085840581
# Row 5
$ws.Range("B5").Value = 0.8964772808081136
$ws.Range("C5").Value = 0.1555753714161483
$ws.Range("E5").Value = 0.09024695334857158
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.8308822803044507
$ws.Range("H5").Value = 0.8904958936483212
$ws.Range("I5").Value = 0.8269908835474808
$ws.Range("L5").Value = 0.2026865156704929
$ws.Range("M5").Value = 0.213786094590759
$ws.Range("N5").Value = 1.441309367947799
# Row 6
$ws.Range("B6").Value = 0.8929389116609343
$ws.Range("C6").Value = 0.1548075608046986
$ws.Range("E6").Value = 0.09026570664038314
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.8307542821724354
$ws.Range("H6").Value = 0.8906908155661029
$ws.Range("I6").Value = 0.8273038757487328
$ws.Range("L6").Value = 0.2025863876267309
$ws.Range("M6").Value = 0.213234395609355
$ws.Range("N6").Value = 1.442083464083321
# Row 7
$ws.Range("B7").Value = 0.9175278370598505
$ws.Range("C7").Value = 0.1601297722981769
$ws.Range("E7").Value = 0.09013690913426664
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.8317039924261564
$ws.Range("H7").Value = 0.8893838096243059
$ws.Range("I7").Value = 0.8251779162322421
$ws.Range("L7").Value = 0.2032908008045382
$ws.Range("M7").Value = 0.2170728131584667
$ws.Range("N7").Value = 1.436761663090547
# Row 8
$ws.Range("B8").Value = 1.026740653957972
$ws.Range("C8").Value = 0.1834368170488574
$ws.Range("E8").Value = 0.08960241767724286
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.837414776237452
$ws.Range("H8").Value = 0.8847546911140398
$ws.Range("I8").Value = 0.8169499191504883
$ws.Range("L8").Value = 0.2066319772872021
$ws.Range("M8").Value = 0.2342334394847398
$ws.Range("N8").Value = 1.414546562587223
# Row 9
$ws.Range("B9").Value = 1.243413859792724
$ws.Range("C9").Value = 0.2285041083789565
$ws.Range("E9").Value = 0.08867543998649463
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.854116368449553
$ws.Range("H9").Value = 0.8797791264002939
$ws.Range("I9").Value = 0.8049757507656921
$ws.Range("L9").Value = 0.2140190348238491
$ws.Range("M9").Value = 0.2686797731844237
$ws.Range("N9").Value = 1.375544119948071
# Row 10
$ws.Range("B10").Value = 1.404123278357019
$ws.Range("C10").Value = 0.2612344658201096
$ws.Range("E10").Value = 0.08806772908958882
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.8697825429484709
$ws.Range("H10").Value = 0.8786356868334053
$ws.Range("I10").Value = 0.798730074179673
$ws.Range("L10").Value = 0.2199551448359784
$ws.Range("M10").Value = 0.2944707039993162
$ws.Range("N10").Value = 1.349669654958333
# Row 11
$ws.Range("B11").Value = 1.477563410513142
$ws.Range("C11").Value = 0.2760436631714924
$ws.Range("E11").Value = 0.08780708011078375
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.8776557957110356
$ws.Range("H11").Value = 0.8786635434740191
$ws.Range("I11").Value = 0.7964461721145781
$ws.Range("L11").Value = 0.2227661450364451
$ws.Range("M11").Value = 0.3063085603922389
$ws.Range("N11").Value = 1.33850253605107
# Row 12
$ws.Range("B12").Value = 1.505420665834436
$ws.Range("C12").Value = 0.2816401050105526
$ws.Range("E12").Value = 0.08771064308650223
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.8807452324366096
$ws.Range("H12").Value = 0.8787530876855243
$ws.Range("I12").Value = 0.7956617220856259
$ws.Range("L12").Value = 0.2238464901631687
$ws.Range("M12").Value = 0.3108063436367203
$ws.Range("N12").Value = 1.334360639074152
# Row 13
$ws.Range("B13").Value = 1.49941902597385
$ws.Range("C13").Value = 0.280435321288337
$ws.Range("E13").Value = 0.08773131189859362
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.8800750530866139
$ws.Range("H13").Value = 0.8787302863154451
$ws.Range("I13").Value = 0.7958270865718902
$ws.Range("L13").Value = 0.2236131126802121
$ws.Range("M13").Value = 0.3098369983101108
$ws.Range("N13").Value = 1.335248806527947
# Row 14
$ws.Range("B14").Value = 1.479854304618414
$ws.Range("C14").Value = 0.2765043157357923
$ws.Range("E14").Value = 0.08779910082388709
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.8779077972291986
$ws.Range("H14").Value = 0.8786693260373681
$ws.Range("I14").Value = 0.7963800217295187
$ws.Range("L14").Value = 0.2228547075769143
$ws.Range("M14").Value = 0.3066782948910358
$ws.Range("N14").Value = 1.338160038562194
# Row 15
$ws.Range("B15").Value = 1.467876455526209
$ws.Range("C15").Value = 0.274094967075257
$ws.Range("E15").Value = 0.08784091826722884
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.8765943751947844
$ws.Range("H15").Value = 0.8786422793778854
$ws.Range("I15").Value = 0.7967291910101366
$ws.Range("L15").Value = 0.2223922298001639
$ws.Range("M15").Value = 0.3047454535474898
$ws.Range("N15").Value = 1.339954565372704
# Row 16
$ws.Range("B16").Value = 1.399330421065997
$ws.Range("C16").Value = 0.2602650421345913
$ws.Range("E16").Value = 0.08808508046931973
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.8692830863448648
$ws.Range("H16").Value = 0.87864491012472
$ws.Range("I16").Value = 0.7988905702353861
$ws.Range("L16").Value = 0.2197736631406286
$ws.Range("M16").Value = 0.2936991837812641
$ws.Range("N16").Value = 1.350411600404243
# Row 17
$ws.Range("B17").Value = 1.357364297692584
$ws.Range("C17").Value = 0.2517603571511415
$ws.Range("E17").Value = 0.08823890814584934
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.8649895339832057
$ws.Range("H17").Value = 0.8787870192936964
$ws.Range("I17").Value = 0.8003594376503287
$ws.Range("L17").Value = 0.2181955728082556
$ws.Range("M17").Value = 0.2869495752226072
$ws.Range("N17").Value = 1.356981238882675
# Row 18
$ws.Range("B18").Value = 1.333257845394655
$ws.Range("C18").Value = 0.2468611551717856
$ws.Range("E18").Value = 0.08832887351904495
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.8625902342076586
$ws.Range("H18").Value = 0.8789203267158427
$ws.Range("I18").Value = 0.8012567292318238
$ws.Range("L18").Value = 0.2172983120748739
$ws.Range("M18").Value = 0.2830773096545727
$ws.Range("N18").Value = 1.36081670112754
# Row 19
$ws.Range("B19").Value = 1.325101222770002
$ws.Range("C19").Value = 0.2452010751621572
$ws.Range("E19").Value = 0.08835958999300664
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.8617899164379565
$ws.Range("H19").Value = 0.8789743136995298
$ws.Range("I19").Value = 0.8015695345937885
$ws.Range("L19").Value = 0.2169963046881378
$ws.Range("M19").Value = 0.2817679362705405
$ws.Range("N19").Value = 1.36212507197714
# Row 20
$ws.Range("B20").Value = 1.361828423872964
$ws.Range("C20").Value = 0.2526664747553298
$ws.Range("E20").Value = 0.08822237900168062
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.8654393159460909
$ws.Range("H20").Value = 0.8787665528411424
$ws.Range("I20").Value = 0.8001976450765085
$ws.Range("L20").Value = 0.2183624855055513
$ws.Range("M20").Value = 0.2876670556322694
$ws.Range("N20").Value = 1.356276012468864
# Row 21
$ws.Range("B21").Value = 1.485599666049666
$ws.Range("C21").Value = 0.2776592587100879
$ws.Range("E21").Value = 0.08777912815660993
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.8785414364104724
$ws.Range("H21").Value = 0.8786850860954871
$ws.Range("I21").Value = 0.7962154266411332
$ws.Range("L21").Value = 0.2230770387780012
$ws.Range("M21").Value = 0.3076056752935372
$ws.Range("N21").Value = 1.337302581323044
# Row 22
$ws.Range("B22").Value = 1.566765346915133
$ws.Range("C22").Value = 0.2939265869183885
$ws.Range("E22").Value = 0.08750263747984621
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.8877342749287038
$ws.Range("H22").Value = 0.8790923794176422
$ws.Range("I22").Value = 0.794081661149292
$ws.Range("L22").Value = 0.2262508217146859
$ws.Range("M22").Value = 0.3207243498344567
$ws.Range("N22").Value = 1.325408614617601
# Row 23
$ws.Range("B23").Value = 1.523420871613723
$ws.Range("C23").Value = 0.2852505206732587
$ws.Range("E23").Value = 0.08764900037045642
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.8827700404994658
$ws.Range("H23").Value = 0.8788327963340663
$ws.Range("I23").Value = 0.795177500827343
$ws.Range("L23").Value = 0.2245484550919628
$ws.Range("M23").Value = 0.3137146862626352
$ws.Range("N23").Value = 1.331710292591207
# Row 24
$ws.Range("B24").Value = 1.359810128728441
$ws.Range("C24").Value = 0.2522568493476456
$ws.Range("E24").Value = 0.08822984706869552
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.8652357542598708
$ws.Range("H24").Value = 0.8787756449790578
$ws.Range("I24").Value = 0.800270626983135
$ws.Range("L24").Value = 0.2182869931111213
$ws.Range("M24").Value = 0.2873426572683684
$ws.Range("N24").Value = 1.356594663137123
# Row 25
$ws.Range("B25").Value = 1.184530541071979
$ws.Range("C25").Value = 0.2163796418379036
$ws.Range("E25").Value = 0.0889132953563746
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.8490048175773666
$ws.Range("H25").Value = 0.8806848220859536
$ws.Range("I25").Value = 0.8077680667160791
$ws.Range("L25").Value = 0.2119312687333164
$ws.Range("M25").Value = 0.2592761967748842
$ws.Range("N25").Value = 1.385607021892106
